$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-18 from
# 2023-10-04 (serial 45203) to 2023-10-06 (serial 45205), keeping the
# existing date formatting/style on the cells untouched.
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
